$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.467.00'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.878.63'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.83'
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4756'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.77%  '
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.82'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07740'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.77'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7382'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +7.89%  '
$ws.Range("D14").Value = '1.878.95'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.127'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '272.98'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").Value = '30.472.85'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007588'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("D21").Value = '2.126.19'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.171'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.296'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.26'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.90'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.943'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.55%  '
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09974'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.516'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.316'
$ws.Range("D32").ClearFormats()
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("E34").Value = '  +2.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.124'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6998'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.740'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.333'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.937'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.06'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4182'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9997'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8392'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.69'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.249'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.081'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.48'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '917.24'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05633'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.01%  '
